# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of cell -> new value to apply to each of the target sheets.
$updates = @{
    "F8"  = 11139
    "F9"  = 4274
    "F11" = 21
    "F15" = 91
    "F18" = 481
    "F19" = 11217
    "F20" = 11061
    "F22" = 35
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
